$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 828228.8
$ws.Range("J17").Value = 978770.4399999999
$ws.Range("L17").Value = 2936311.32
$ws.Range("N17").Value = -2936647.32

# Hunk 1: sheet ALC row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 3839
$ws.Range("I31").Value = 3839
$ws.Range("K31").Value = 11517
$ws.Range("M31").Value = -11287

# Hunk 2: sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 140
$ws.Range("I33").Value = 136.42857
$ws.Range("J33").Value = 152.5
$ws.Range("K33").Value = 136.42857
$ws.Range("L33").Value = 152.5
$ws.Range("M33").Value = 92.57142999999999
$ws.Range("N33").Value = -610.5

# Hunk 3: sheet ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6819219.5
$ws.Range("J112").Value = 7793296.5
$ws.Range("L112").Value = 23379889.5
$ws.Range("N112").Value = -23382105.5

# Hunk 4: sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5711598
$ws.Range("I138").Value = 1062448.8
$ws.Range("J138").Value = 13336203
$ws.Range("K138").Value = 3187346.4
$ws.Range("L138").Value = 40008609
$ws.Range("M138").Value = -3182206.4
$ws.Range("N138").Value = -40018889

# Hunk 5: sheet ARM row 43
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 6000
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6626

# Hunk 6: sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2102.4656
$ws.Range("I132").Value = 1954.2821
$ws.Range("J132").Value = 2406.6316
$ws.Range("K132").Value = 5862.846299999999
$ws.Range("L132").Value = 7219.8948
$ws.Range("M132").Value = -3332.846299999999
$ws.Range("N132").Value = -12279.8948

# Hunk 7: sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1743.2084
$ws.Range("I58").Value = 930.4706
$ws.Range("K58").Value = 930.4706
$ws.Range("M58").Value = -727.4706

# Hunk 8: sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5210141
$ws.Range("I99").Value = 8929854
$ws.Range("J99").Value = 2542.8
$ws.Range("K99").Value = 8929854
$ws.Range("L99").Value = 2542.8
$ws.Range("M99").Value = -8928356
$ws.Range("N99").Value = -5538.8

# Hunk 9: sheet CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 321.8421
$ws.Range("I107").Value = 201.07143
$ws.Range("J107").Value = 660
$ws.Range("K107").Value = 201.07143
$ws.Range("L107").Value = 660
$ws.Range("M107").Value = 1718.92857
$ws.Range("N107").Value = -4500

# Hunk 10: sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5210141
$ws.Range("I126").Value = 8929854
$ws.Range("J126").Value = 2542.8
$ws.Range("K126").Value = 26789562
$ws.Range("L126").Value = 7628.400000000001
$ws.Range("M126").Value = -26787092
$ws.Range("N126").Value = -12568.4

# Hunk 11: sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1459.1837
$ws.Range("I132").Value = 1406.2709
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 4218.8127
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -1688.8127
$ws.Range("N132").Value = -17057

# Hunk 12: sheet CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 21148.428
$ws.Range("J133").Value = 22715.143
$ws.Range("L133").Value = 22715.143
$ws.Range("N133").Value = -27775.143

# Hunk 13: sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1743.2084
$ws.Range("I136").Value = 930.4706
$ws.Range("K136").Value = 2791.4118
$ws.Range("M136").Value = -241.4117999999999

# Hunk 14: sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1226.4222
$ws.Range("I5").Value = 870.65515
$ws.Range("J5").Value = 1871.25
$ws.Range("K5").Value = 2611.96545
$ws.Range("L5").Value = 5613.75
$ws.Range("M5").Value = -2499.96545
$ws.Range("N5").Value = -5837.75

# Hunk 15: sheet CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 998
$ws.Range("J75").Value = 997.6667
$ws.Range("L75").Value = 2993.0001
$ws.Range("N75").Value = -4989.0001

# Hunk 16: sheet CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 998
$ws.Range("J78").Value = 997.6667
$ws.Range("L78").Value = 8979.0003
$ws.Range("N78").Value = -18963.0003

# Hunk 17: sheet CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 602127.6
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 751909.5
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 2255728.5
$ws.Range("M103").Value = -8121
$ws.Range("N103").Value = -2257486.5

# Hunk 18: sheet CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1091.3334
$ws.Range("J117").Value = 1527
$ws.Range("L117").Value = 4581
$ws.Range("N117").Value = -11465

# Hunk 19: sheet CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1674.3334
$ws.Range("J132").Value = 1674.3334
$ws.Range("L132").Value = 15069.0006
$ws.Range("N132").Value = -20129.0006

# Hunk 20: sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1226.4222
$ws.Range("I135").Value = 870.65515
$ws.Range("J135").Value = 1871.25
$ws.Range("K135").Value = 7835.896350000001
$ws.Range("L135").Value = 16841.25
$ws.Range("M135").Value = -5300.896350000001
$ws.Range("N135").Value = -21911.25

# Hunk 21: sheet CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4594773.5
$ws.Range("I137").Value = 7146101
$ws.Range("J137").Value = 129950.25
$ws.Range("K137").Value = 21438303
$ws.Range("L137").Value = 389850.75
$ws.Range("M137").Value = -21433203
$ws.Range("N137").Value = -400050.75

# Hunk 22: sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 556731.8
$ws.Range("I122").Value = 618340.4399999999
$ws.Range("J122").Value = 2254
$ws.Range("K122").Value = 1855021.32
$ws.Range("L122").Value = 6762
$ws.Range("M122").Value = -1852571.32
$ws.Range("N122").Value = -11662

# Hunk 23: sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8712.666999999999
$ws.Range("I22").Value = 1283.3334
$ws.Range("J22").Value = 13665.556
$ws.Range("K22").Value = 1283.3334
$ws.Range("L22").Value = 13665.556
$ws.Range("M22").Value = -988.3334
$ws.Range("N22").Value = -14255.556

# Hunk 24: sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8712.666999999999
$ws.Range("I27").Value = 1283.3334
$ws.Range("J27").Value = 13665.556
$ws.Range("K27").Value = 1283.3334
$ws.Range("L27").Value = 13665.556
$ws.Range("M27").Value = -1176.3334
$ws.Range("N27").Value = -13879.556

# Hunk 25: sheet LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5433.725
$ws.Range("I136").Value = 3877.68
$ws.Range("K136").Value = 11633.04
$ws.Range("M136").Value = -9083.039999999999

# Hunk 26: sheet WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 11000
$ws.Range("I32").Value = 11000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -10683
$ws.Range("N32").ClearContents()

# Hunk 27: sheet WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7400
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Hunk 28: sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 585.25
$ws.Range("I113").Value = 538.8570999999999
$ws.Range("J113").Value = 621.3333
$ws.Range("K113").Value = 1616.5713
$ws.Range("L113").Value = 1863.9999
$ws.Range("M113").Value = 553.4287000000002
$ws.Range("N113").Value = -6203.9999

# Hunk 29: sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 85491.836
$ws.Range("I122").Value = 126737.75
$ws.Range("K122").Value = 380213.25
$ws.Range("M122").Value = -377763.25

# Hunk 30: sheet WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 25889.092
$ws.Range("J123").Value = 25889.092
$ws.Range("L123").Value = 25889.092
$ws.Range("N123").Value = -35689.092

# Hunk 31: sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9617469
$ws.Range("I132").Value = 11365552
$ws.Range("J132").Value = 3015
$ws.Range("K132").Value = 34096656
$ws.Range("L132").Value = 9045
$ws.Range("M132").Value = -34094126
$ws.Range("N132").Value = -14105

# Hunk 32: sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19106.861
$ws.Range("I136").Value = 25697.928
$ws.Range("K136").Value = 77093.784
$ws.Range("M136").Value = -74543.784
